# Applies updated cryptocurrency price/volume data to Sheet1 (cryptos.xlsx),
# matching the "Updated cryptos list" GitHub Actions commit.
#
# Price values in column D are mostly plain decimal numbers (e.g. "617.92").
# Excel would silently reinterpret those as numeric values (and normalize away
# trailing zeros / leading zeros, e.g. "1.00" -> 1), so for the cells where the
# new price text would otherwise be parsed as a number we prefix the value with
# a leading apostrophe (Excel's standard "treat as text" marker). The apostrophe
# itself is not stored in the cell - Excel only uses it to force text entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = [char]39   # single-quote / apostrophe, used as the "force text" prefix

# Row 2
$ws.Range("D2").Value = '69.608.84'
$ws.Range("E2").Value = '  -1.85%  '

# Row 3
$ws.Range("D3").Value = '3.750.57'
$ws.Range("E3").Value = '  +2.55%  '

# Row 4
$ws.Range("D4").Value = "$q" + '0.999'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = "$q" + '617.92'
$ws.Range("E5").Value = '  +3.65%  '

# Row 6
$ws.Range("D6").Value = "$q" + '177.34'
$ws.Range("E6").Value = '  -2.81%  '

# Row 7
$ws.Range("D7").Value = '3.748.08'
$ws.Range("E7").Value = '  +2.63%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = "$q" + '0.532'
$ws.Range("E9").Value = '  -0.37%  '

# Row 10
$ws.Range("D10").Value = "$q" + '0.168'
$ws.Range("E10").Value = '  +3.71%  '

# Row 11
$ws.Range("D11").Value = "$q" + '6.27'
$ws.Range("E11").Value = '  -4.52%  '

# Row 12
$ws.Range("D12").Value = "$q" + '0.486'
$ws.Range("E12").Value = '  -2.03%  '

# Row 13
$ws.Range("D13").Value = "$q" + '40.77'
$ws.Range("E13").Value = '  +1.19%  '

# Row 14
$ws.Range("E14").Value = '  +1.55%  '

# Row 15
$ws.Range("D15").Value = '4.376.39'
$ws.Range("E15").Value = '  +2.48%  '

# Row 16
$ws.Range("D16").Value = '3.746.62'
$ws.Range("E16").Value = '  +2.24%  '

# Row 17
$ws.Range("D17").Value = '69.556.62'
$ws.Range("E17").Value = '  -1.93%  '

# Row 18
$ws.Range("E18").Value = '  +0.16%  '

# Row 19
$ws.Range("D19").Value = "$q" + '7.53'
$ws.Range("E19").Value = '  +0.94%  '

# Row 20
$ws.Range("E20").Value = '  -1.86%  '

# Row 21
$ws.Range("D21").Value = "$q" + '507.94'
$ws.Range("E21").Value = '  -0.55%  '

# Row 22
$ws.Range("D22").Value = "$q" + '9.51'
$ws.Range("E22").Value = '  +4.09%  '

# Row 23
$ws.Range("D23").Value = "$q" + '0.721'
$ws.Range("E23").Value = '  -2.36%  '

# Row 24
$ws.Range("D24").Value = "$q" + '2.51'
$ws.Range("E24").Value = '  +2.20%  '

# Row 25
$ws.Range("D25").Value = "$q" + '86.59'
$ws.Range("E25").Value = '  -0.81%  '

# Row 26
$ws.Range("D26").Value = "$q" + '13.13'
$ws.Range("E26").Value = '  -2.40%  '

# Row 27
$ws.Range("D27").Value = "$q" + '11.08'
$ws.Range("E27").Value = '  +1.69%  '

# Row 28
$ws.Range("E28").Value = '  +23.32%  '

# Row 29
$ws.Range("E29").Value = '  -0.19%  '

# Row 30
$ws.Range("D30").Value = "$q" + '2.48'
$ws.Range("E30").Value = '  -1.65%  '

# Row 31
$ws.Range("E31").Value = '  +4.79%  '

# Row 32
$ws.Range("D32").Value = "$q" + '7.80'
$ws.Range("E32").Value = '  -4.22%  '

# Row 33
$ws.Range("D33").Value = "$q" + '30.94'
$ws.Range("E33").Value = '  -1.46%  '

# Row 34
$ws.Range("E34").Value = '  -1.57%  '

# Row 35
$ws.Range("D35").Value = "$q" + '0.997'
$ws.Range("E35").Value = '  -0.24%  '

# Row 36
$ws.Range("E36").Value = '  +5.09%  '

# Row 37
$ws.Range("D37").Value = "$q" + '6.13'
$ws.Range("E37").Value = '  +0.92%  '

# Row 38
$ws.Range("D38").Value = "$q" + '0.334'
$ws.Range("E38").Value = '  -3.05%  '

# Row 39
$ws.Range("E39").Value = '  +1.73%  '

# Row 40
$ws.Range("E40").Value = '  -1.73%  '

# Row 41
$ws.Range("D41").Value = "$q" + '50.51'
$ws.Range("E41").Value = '  -0.99%  '

# Row 42
$ws.Range("D42").Value = "$q" + '44.82'
$ws.Range("E42").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Value = "$q" + '420.62'
$ws.Range("E43").Value = '  +1.92%  '

# Row 44
$ws.Range("D44").Value = "$q" + '8.66'
$ws.Range("E44").Value = '  -1.61%  '

# Row 45
$ws.Range("D45").Value = '3.008.12'
$ws.Range("E45").Value = '  -3.83%  '

# Row 46
$ws.Range("E46").Value = '  -0.56%  '

# Row 47
$ws.Range("D47").Value = "$q" + '0.0360'
$ws.Range("E47").Value = '  -2.00%  '

# Rows 48-50: Monero overtook InjectiveProtocol and USDe in rank, so the three
# coins shift position (Monero: 50 -> 48, InjectiveProtocol: 48 -> 49, USDe: 49 -> 50)
# while also getting refreshed price/volume figures.
# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = "$q" + '139.14'
$ws.Range("E48").Value = '  +1.21%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "$q" + '27.25'
$ws.Range("E49").Value = '  -3.21%  '

# Row 50
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").Value = "$q" + '1.00'
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("D51").Value = "$q" + '2.48'
$ws.Range("E51").Value = '  +0.39%  '
